$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 -- this shifts rows 3..19 (A,B,C,D columns, and existing E values) up by one,
# so that new row N (2..18) ends up holding what used to be row N+1's A/B/C/D data,
# and the oldest data point (previously row 2) is dropped.
$ws.Rows(2).Delete()

# Now overwrite column E (y_1_forecast) with the freshly computed forecast values,
# including the previously-empty E2.
$eValues = @(
  "0.8212989654785341",
  "1.183007486132071",
  "1.015842920196763",
  "0.9092565586104273",
  "1.236730309040235",
  "1.029015928490629",
  "1.358148715145191",
  "1.528208222695326",
  "1.634928000057778",
  "1.67176973076042",
  "1.603287858019664",
  "0.8408455317168162",
  "-1.875058665585216",
  "5.03478667886097",
  "2.399708479013141",
  "0.8520283695166997",
  "0.299857156820571"
)

for ($i = 0; $i -lt $eValues.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 5).Value = [double]$eValues[$i]
}
